# Reproduce the commit: "Made sheet2 the same as reader/sheet2 ..."
# Content changes:
#  - Sheet2 ("Sheet2 - Numbers") gains a new column AA (rows 1-30) with
#    values 100..129, becomes the active/selected sheet/tab, and its
#    selection moves to AA1:AA30 (active cell AA1).
#  - Sheet4 ("Sheet4 - Dates") loses the tabSelected flag (because Sheet2
#    becomes the active tab instead) and its PageSetup paper size changes
#    from "use printer default" (0) to A4 (9).
#  - Workbook style names / default pivot style are locale/version
#    artifacts of the original authoring app and are attempted here on a
#    best-effort basis.

$wb = $excel.ActiveWorkbook

# ---- Sheet2: add column AA with values 100-129 in rows 1-30 ----
$ws2 = $wb.Worksheets.Item("Sheet2 - Numbers")
for ($row = 1; $row -le 30; $row++) {
    $ws2.Cells.Item($row, 27).Value = 99 + $row
}

# Make Sheet2 the active sheet/tab (activeTab moves from Sheet4 to Sheet2,
# and tabSelected moves with it).
$ws2.Activate()

# Scroll so column O is the left-most visible column, then select AA1:AA30
# with AA1 as the active cell (matches topLeftCell="O1" / selection in the
# target worksheet XML as closely as this host's COM surface allows).
$excel.ActiveWindow.ScrollColumn = 15
$ws2.Range("AA1:AA30").Select() | Out-Null

# ---- Sheet4: page setup paper size 0 (default) -> 9 (A4) ----
$ws4 = $wb.Worksheets.Item("Sheet4 - Dates")
$ws4.PageSetup.PaperSize = 9

# ---- Workbook style names (locale rename Pourcentage/Normal -> Prozent/Standard) ----
$styles = $wb.Styles
for ($i = 1; $i -le $styles.Count; $i++) {
    $style = $styles.Item($i)
    if ($style.Name -eq "Normal") {
        $style.Name = "Standard"
    } elseif ($style.Name -eq "Pourcentage") {
        $style.Name = "Prozent"
    }
}

# ---- Default pivot table style ----
$wb.DefaultPivotTableStyle = "PivotStyleMedium4"
